# Update TSRs in ENA templates
$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump version number ---
$wsMeta = $wb.Worksheets.Item("isa_template")
$wsMeta.Range("B4").Value = "1.0.2"

# --- New Table sheet: update column headers + TSR/TAN data ---
$wsTable = $wb.Worksheets.Item("New Table")

# Header renames
$wsTable.Range("A1").Value = "Input [Data]"
$wsTable.Range("B1").Value = "Characteristic [organism]"
$wsTable.Range("C1").Value = "Term Source REF (OBI:0100026)"
$wsTable.Range("D1").Value = "Term Accession Number (OBI:0100026)"
$wsTable.Range("Z1").Value = "Output [Data]"

# Data row updates
$wsTable.Range("E2").Value = "No"
$wsTable.Range("F2").Value = "NCIT"
$wsTable.Range("G2").Value = "https://bioregistry.io/NCIT:C49487"
$wsTable.Range("N2").Value = "No"
$wsTable.Range("O2").Value = "NCIT"
$wsTable.Range("P2").Value = "https://bioregistry.io/NCIT:C49487"
$wsTable.Range("Q2").Value = "No"
$wsTable.Range("R2").Value = "NCIT"
$wsTable.Range("S2").Value = "https://bioregistry.io/NCIT:C49487"
